$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.409.03'
$ws.Range("E2").Value = '  -0.70%  '

$ws.Range("D3").Value = '3.779.23'
$ws.Range("E3").Value = '  -1.56%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '647.67'
$ws.Range("E5").Value = '  +1.66%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.05'
$ws.Range("E6").Value = '  -0.43%  '

$ws.Range("D7").Value = '3.777.26'
$ws.Range("E7").Value = '  -1.56%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("E9").Value = '  +0.69%  '

$ws.Range("E10").Value = '  -2.35%  '

$ws.Range("E11").Value = '  +0.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.91'
$ws.Range("E12").Value = '  +3.56%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000240'
$ws.Range("E13").Value = '  -4.99%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.96'
$ws.Range("E14").Value = '  -3.08%  '

$ws.Range("D15").Value = '4.408.13'
$ws.Range("E15").Value = '  -1.63%  '

$ws.Range("D16").Value = '3.785.67'
$ws.Range("E16").Value = '  -4.98%  '

$ws.Range("D17").Value = '69.313.92'
$ws.Range("E17").Value = '  -0.70%  '

$ws.Range("E18").Value = '  -1.95%  '

$ws.Range("E19").Value = '  +0.21%  '

$ws.Range("E20").Value = '  -1.87%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '468.03'
$ws.Range("E21").Value = '  -0.16%  '

$ws.Range("E22").Value = '  -1.56%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.709'
$ws.Range("E23").Value = '  -0.10%  '

$ws.Range("E24").Value = '  -5.66%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.86'
$ws.Range("E25").Value = '  -2.37%  '

$ws.Range("E26").Value = '  +2.52%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.40'
$ws.Range("E27").Value = '  +2.93%  '

$ws.Range("E28").Value = '  -3.50%  '

$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("D30").Value = '3.926.42'
$ws.Range("E30").Value = '  -1.50%  '

$ws.Range("E31").Value = '  -0.23%  '

$ws.Range("E32").Value = '  +2.20%  '

$ws.Range("E33").Value = '  -2.33%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '28.74'
$ws.Range("E34").Value = '  -2.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.173'
$ws.Range("E35").Value = '  +14.64%  '

$ws.Range("E36").Value = '  -0.02%  '

$ws.Range("D37").Value = '3.732.72'
$ws.Range("E37").Value = '  -1.30%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.87'
$ws.Range("E38").Value = '  -2.36%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.101'
$ws.Range("E39").Value = '  -2.44%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.85'
$ws.Range("E40").Value = '  -1.42%  '

$ws.Range("E41").Value = '  -6.31%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  +0.00%  '

$ws.Range("E43").Value = '  -2.63%  '

$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '45.39'
$ws.Range("E45").Value = '  +3.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.00'
$ws.Range("E46").Value = '  +2.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '157.09'
$ws.Range("E47").Value = '  -0.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '47.55'
$ws.Range("E48").Value = '  +0.49%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.298'
$ws.Range("E49").Value = '  -1.46%  '

$ws.Range("E50").Value = '  -0.94%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.39'
$ws.Range("E51").Value = '  -1.02%  '
